$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Aciers")
$ws.Range("AA1:AB24").ClearContents()
$ws.Range("AA1:AB24").ClearFormats()
Write-Host "UsedRange:" $ws.UsedRange.Address()
